$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.133.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.32%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.877.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.65%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.29%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.22%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.31%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5063"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3842"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.39%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08634"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.34%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.118"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.07%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.58"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.69%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.328"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.31%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.96%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.885.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.56%  "

# Row 15
$ws.Range("B15").Value = "BinanceUSD"
$ws.Range("C15").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.005"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.35%  "

# Row 16
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.180"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.73%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001102"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.98%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.64%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06624"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.17%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.93%  "

# Row 21
$ws.Range("E21").Value = "  +0.09%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.095"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.84%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.154.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.44%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.08%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.269"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.30%  "

# Row 26
$ws.Range("B26").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C26").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.099.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.59%  "

# Row 27
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.549"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.23%  "

# Row 28
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.46%  "

# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.62%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.27%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1052"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.05%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.059"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.93%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.589"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.00%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.601"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.38%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.666"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.17%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02446"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.05%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06543"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.81%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2172"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.73%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.205"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.21%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.241"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.29%  "

# Row 41
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.45%  "

# Row 42
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6372"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.09%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.901"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.16%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.43%  "

# Row 45
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5988"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.52%  "

# Row 46
$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.278"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.05%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.674"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.29%  "

# Row 48
$ws.Range("B48").Value = "EOS"
$ws.Range("C48").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.231"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.71%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.990"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.47%  "

# Row 50
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "121.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.60%  "

# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.85"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.85%  "
